$d = $word.ActiveDocument

# Merge the split runs in the Title paragraph into a single run of text.
$d.Content.Find.Execute(
    "Answers: Trigonometry (radians)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Answers: Trigonometry (radians)", 2) | Out-Null

# Merge the split runs in the Abstract paragraph into a single run of text.
$d.Content.Find.Execute(
    "Answers to the questions on trigonometry, using radians to measure angles.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Answers to the questions on trigonometry, using radians to measure angles.", 2) | Out-Null
